# Generate Report for Handback
# This script brings the localization-status workbook up to date with the
# latest handback run:
#   - The overall status moves from "Ready for handoff" to
#     "Handed back: in sync with en-US" (Overview + per-locale sheets).
#   - The per-locale "Latest Handback DateTime" is refreshed.
#   - The (now resolved) handback-version error on each locale sheet is
#     cleared.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: both locale-status columns reflect the new status text.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E1").ColumnWidth = 29.166666666666668
$wsOverview.Range("F1").ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# zh-cn sheet: status, refreshed handback time, cleared error detail.
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("K2").Value = "2016-08-14 03:09:01"
$wsZhCn.Range("P2").Value = ""
$wsZhCn.Range("C1").ColumnWidth = 29.166666666666668
$wsZhCn.Range("P1").ColumnWidth = 12.833333333333332

# ---------------------------------------------------------------------
# de-de sheet: status, refreshed handback time, cleared error detail.
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("K2").Value = "2016-08-14 03:09:13"
$wsDeDe.Range("P2").Value = ""
$wsDeDe.Range("C1").ColumnWidth = 29.166666666666668
$wsDeDe.Range("P1").ColumnWidth = 12.833333333333332
